$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update 想去人数 (F column) for several rows
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 299
$ws1.Range("F4").Value = 8007
$ws1.Range("F5").Value = 5841
$ws1.Range("F6").Value = 496
$ws1.Range("F10").Value = 284
$ws1.Range("F11").Value = 371

# Sheet "全部类型" (sheet4): same underlying data, different row numbers
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 299
$ws4.Range("F4").Value = 8007
$ws4.Range("F5").Value = 5841
$ws4.Range("F6").Value = 496
$ws4.Range("F10").Value = 284
$ws4.Range("F14").Value = 371
